# Apply the "adding averages and more checks" update:
#  - Training Dashboard (sheet 1): refresh "PERIOD TO EXPIRE" (H) and
#    "LAST UPDATE" (I) columns for rows 3-25 to reflect a later
#    last-update date (16-Sep-2025 instead of 08-Sep-2025).
#  - Exam Dashboard (sheet 2): widen the COMMENTS column and replace the
#    generic "OK" comment with a more descriptive "date is valid" for
#    rows 3-9.
#  - Header rows on both sheets get an explicit white, bold font so the
#    text is legible against the dark-blue fill.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Training Dashboard
$ws2 = $wb.Worksheets.Item(2)   # Exam Dashboard

# ---------------------------------------------------------------------
# Training Dashboard: updated "PERIOD TO EXPIRE" / "LAST UPDATE" values
# ---------------------------------------------------------------------
$newPeriod = @{
    3  = 210
    4  = 384
    5  = 321
    6  = 244
    7  = 358
    8  = 364
    9  = 672
    10 = 268
    11 = 245
    12 = 365
    13 = 247
    14 = 355
    15 = 307
    16 = 379
    17 = 392
    18 = 380
    19 = 336
    20 = 85
    21 = -126
    22 = 174
    23 = 177
    24 = 189
    25 = 232
}

for ($row = 3; $row -le 25; $row++) {
    $ws1.Cells.Item($row, 8).Value = $newPeriod[$row]
    # Leading apostrophe forces this date-looking value to stay plain
    # text (matching the rest of the sheet, which stores dates as text).
    $ws1.Cells.Item($row, 9).Value = "'16-Sep-2025"
}

# ---------------------------------------------------------------------
# Exam Dashboard: wider COMMENTS column + more descriptive comment text
# ---------------------------------------------------------------------
$ws2.Columns.Item(5).ColumnWidth = 14.17

for ($row = 3; $row -le 9; $row++) {
    $ws2.Cells.Item($row, 5).Value = "date is valid"
}

# ---------------------------------------------------------------------
# Header / title styling: make the bold header font explicitly white
# ---------------------------------------------------------------------
$ws1.Range("A2:K2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215

$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Color = 16777215
$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215
